# 自动更新价格数据: insert a new top data row (2026-01-10) by pushing
# existing rows down one position (row 2 -> row 3 ... row 51 -> row 52).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 2 (the first data row), shifting
# all existing data rows down by one.
$ws.Rows.Item(2).Insert()

# Row.Insert() copies formatting from the row above (the bold header),
# so reset the new row back to the plain/default style used by every
# other data row before writing values into it.
$newRow = $ws.Range("A2:D2")
$newRow.Style = "Normal"

# Populate the newly inserted row with the latest price data.
# The date is entered with a leading apostrophe so Excel keeps it as
# literal text instead of auto-converting it to a date serial; the
# style is then reset back to the workbook default (matching the
# unformatted data rows) so no stray number format / quote-prefix
# flag is left behind.
$ws.Cells.Item(2, 1).Value = "'2026-01-10"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
